# "Generate Report for Archive" — refresh the localization-status report:
#   - the zh-cn / de-de items have moved from hand-off into translation,
#     so the Status text changes everywhere it is shown
#   - the Status column is narrower now ("In Translation" renders shorter
#     than "Ready for handoff"), so the report's column widths shrink to match

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newWidth  = 12.5   # -> stored column width ~13.33, matching the narrower auto-fit

# --- Overview sheet: Status is shown per-locale in columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet: Status is column C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet: Status is column C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
